$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)
$shape.Table.ApplyStyle("{277EBAE9-94F0-4B25-8C53-B326EDFCCC64}")
